# Insert a new data row at row 133 (pushing the existing rows 133-168 down
# to 134-169) and populate it with the new "Haba" price observation dated
# 2021-11-08 (serial 44508), matching the author's weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(133).Insert()

$ws.Cells.Item(133, 1).Value = 9
$ws.Cells.Item(133, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(133, 3).Value = "Metropolitana"
$ws.Cells.Item(133, 4).Value = 44508
$ws.Cells.Item(133, 5).Value = 13
$ws.Cells.Item(133, 6).Value = 100112026
$ws.Cells.Item(133, 7).Value = "Haba"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 52
$ws.Cells.Item(133, 11).Value = 8000
$ws.Cells.Item(133, 12).Value = 9000
$ws.Cells.Item(133, 13).Value = 8500
$ws.Cells.Item(133, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(133, 15).Value = "Región Metropolitana"
$ws.Cells.Item(133, 16).Value = 340
$ws.Cells.Item(133, 17).Value = 25
$ws.Cells.Item(133, 18).Value = "Hortaliza"
